# Insert a new transaction row above row 2, pushing all existing
# transaction rows down by one (row 2 -> row 3, ..., row 136 -> row 137).
#
# We do this with a bottom-up cascading copy of the four data columns
# (E = Transaction Type, N = Payment Type, P = InternalComment,
#  T = USD Amount) instead of Rows.Insert(), because Insert() drags
# along formatting (it clones the style of the row above into every
# column, generating new style/font records) which is not what
# happened in the source edit - only the data shifted down, the
# original per-column formatting stayed untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 136

for ($r = $lastRow; $r -ge 2; $r--) {
    $dest = $r + 1

    $ws.Range("E" + $dest).Value2 = $ws.Range("E" + $r).Value2
    $ws.Range("N" + $dest).Value2 = $ws.Range("N" + $r).Value2
    $ws.Range("P" + $dest).Value2 = $ws.Range("P" + $r).Value2
    $ws.Range("T" + $dest).Value2 = $ws.Range("T" + $r).Value2
}

# Row 125 now holds what used to be row 124's data, which (unlike the
# rows below it) also carries empty date-formatted placeholder cells
# in K/S/AB. Reproduce that formatting on the now-shifted row.
$dateFmt = "yyyy\-mm\-dd\ hh:mm:ss"
$ws.Range("K125").NumberFormat = $dateFmt
$ws.Range("S125").NumberFormat = $dateFmt
$ws.Range("AB125").NumberFormat = $dateFmt

# Populate the brand-new row 2 with the newly added transaction.
$ws.Range("E2").Value2 = "Withdrawal"
$ws.Range("N2").Value2 = "Credit Card"
$ws.Range("P2").Value2 = "Tradeprof"
$ws.Range("T2").Value2 = 269.235

# Update the selection to reflect where the user ended up after the edit.
$ws.Range("E2:E3").Select()
